# "some sdditional codes on submita free ad link"
#
# Adds a new Test Scenario / Test Case row (TS_002 / TC-002 - "Submit a
# Free Ad" link in the header section) to the "Scenario" sheet, and leaves
# the workbook with the "Scenario" tab active (selection on I6) instead of
# the "test case" tab (selection moves to B9 there).

$wb = $excel.ActiveWorkbook

$wsScenario = $wb.Worksheets.Item("Scenario")
$wsTestCase = $wb.Worksheets.Item("test case")

# --- move the selection on the "test case" sheet (it is losing focus) ---
$wsTestCase.Range("B9").Select()

# --- new row of scenario/test-case data on the "Scenario" sheet ---
$wsScenario.Range("A3").Value = "TS_002"
$wsScenario.Range("B3").Value = "Homepage"
$wsScenario.Range("C3").Value = "Header Section"
$wsScenario.Range("D3").Value = "TS_HP_HS_Links"
$wsScenario.Range("E3").Value = "To check the links in header section"
$wsScenario.Range("G3").Value = "TC-002"
$wsScenario.Range("H3").Value = "TC_HP_HS_AdLink"
$wsScenario.Range("I3").Value = "To check Submik a Free add link"

# Match the look of the row above (A2:B2) for the new A3:B3 cells.
$wsScenario.Range("A2:B2").Copy()
$wsScenario.Range("A3:B3").PasteSpecial(-4122)

# Description-style cells wrap their text, same as elsewhere in the sheet.
$wsScenario.Range("E3").WrapText = $true
$wsScenario.Range("G3").WrapText = $true
$wsScenario.Range("I3").WrapText = $true

# --- "Scenario" tab becomes the active tab, cursor parked on I6 ---
$wsScenario.Activate()
$wsScenario.Range("I6").Select()
